$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ARKCORR-22: Minor changes to the dueDate setting for the 'Release' queue.
# Row 31 is the "Set Due Date Release Queue" rule:
#   CONDITION (C31): drop the redundant "dueDate == null &&" guard
#   ACTION    (D31): reset the due date to null instead of "now"
$ws.Range("C31").Value = 'queue.name == "Release"'
$ws.Range("D31").Value = "setDueDate, null"

# Keep the sheet's selection/scroll state pointing at the row that changed.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 19
[void]$ws.Range("D31").Select()
